$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.930.57'
$ws.Range('E2').Value = '  +2.61%  '
$ws.Range('D3').Value = '3.450.32'
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '189.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.632'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '3.442.54'
$ws.Range('E8').Value = '  +1.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.646'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.15'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000277'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('D15').Value = '4.005.45'
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.93'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('D17').Value = '3.456.87'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '66.990.67'
$ws.Range('E18').Value = '  +2.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '480.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.26%  '
$ws.Range('E23').Value = '  +9.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +18.12%  '
$ws.Range('E25').Value = '  +6.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('E27').Value = '  +3.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.95'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.08%  '
$ws.Range('E29').Value = '  +3.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('E31').Value = '  +14.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '600.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.24'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.84%  '
$ws.Range('E35').Value = '  +4.02%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.148'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.391'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.50'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.85%  '
$ws.Range('D41').Value = '0.0₃0754'
$ws.Range('E41').Value = '  +1.83%  '
$ws.Range('D42').Value = '3.216.61'
$ws.Range('E42').Value = '  +4.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0431'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.47%  '
$ws.Range('E45').Value = '  +5.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.79'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +23.81%  '
$ws.Range('E47').Value = '  +2.01%  '
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.74%  '
$ws.Range('E51').Value = '  +6.43%  '
